# Change the table style (Table Design gallery) applied to the cash-flow
# glossary table on slide 16 from "Table_0" ({28B35176-89C3-4D0D-ACE2-188A2F0F83E7})
# to the built-in style {C4574F6F-1FA0-49DB-B79B-1AD107C562D2}.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{C4574F6F-1FA0-49DB-B79B-1AD107C562D2}")
    }
}
